$d = $word.ActiveDocument

# Locate the paragraph holding the "Iso Burkina Faso" label (it also
# contains the _GoBack bookmark).
$isoLabelIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Iso Burkina Faso*") {
        $isoLabelIndex = $i
        break
    }
}

$labelPara = $d.Paragraphs.Item($isoLabelIndex)

# Remove everything after this paragraph's mark through to the very end of
# the document in one shot: this drops the ISO url paragraph and the two
# trailing empty paragraphs (and, since it spans their paragraph marks too,
# merges them away instead of leaving empty paragraphs behind).
$tailRange = $d.Range($labelPara.Range.End, $d.Content.End)
if ($tailRange.Start -lt $tailRange.End) {
    $tailRange.Delete()
}

# Remove the "Iso Burkina Faso" text itself (but not the paragraph mark,
# which still carries the _GoBack bookmark).
$labelPara = $d.Paragraphs.Item($isoLabelIndex)
$textRange = $d.Range($labelPara.Range.Start, $labelPara.Range.End - 1)
if ($textRange.Start -lt $textRange.End) {
    $textRange.Delete()
}

# Drop the paragraph mark's stored run formatting (rStyle/color/underline)
# by restoring the paragraph to the default "Normal" style, leaving a bare
# <w:p> that only still carries the bookmark.
$labelPara = $d.Paragraphs.Item($isoLabelIndex)
$labelPara.Style = "Normal"
